$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-5 from 45212 to 45221
$ws.Range("C2").Value = 45221
$ws.Range("C3").Value = 45221
$ws.Range("C4").Value = 45221
$ws.Range("C5").Value = 45221
